# Update title and a couple of data points per commit "Updated slides and sample"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Visual Performance")

# Title cell B1: "HttpClient Performance" -> "Controls Performance"
$ws.Range("B1").Value = "Controls Performance"

# Data updates in column D (100 series)
$ws.Range("D5").Value = 198
$ws.Range("D7").Value = 204

# Update selection to match the recorded view state
$ws.Range("B15").Select()
